$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.802.79"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "2.218.69"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.24"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.00"
$ws.Range("E6").Value = "  +4.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.82"
$ws.Range("E10").Value = "  +6.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0785"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.16"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.35"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "2.564.65"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.04"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "2.211.77"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.730"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").Value = "39.808.14"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "0.0₃0883"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.04"
$ws.Range("E21").Value = "  +7.69%  "
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.54"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.15"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.46"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.83"
$ws.Range("E27").Value = "  +2.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.72"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("E29").Value = "  +4.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.24"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.87"
$ws.Range("E31").Value = "  +4.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.73"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.93"
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0718"
$ws.Range("E35").Value = "  +4.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("E37").Value = "  +7.32%  "
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0994"
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.77"
$ws.Range("E40").Value = "  +4.61%  "
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("D43").Value = "2.064.68"
$ws.Range("E43").Value = "  +9.07%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.90"
$ws.Range("E44").Value = "  +12.21%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0268"
$ws.Range("E45").Value = "  +4.14%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.10"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.94"
$ws.Range("E47").Value = "  +10.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.61"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "2.435.70"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.54"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.04"
$ws.Range("E51").Value = "  +2.78%  "
